# Delete the "pera" row (row 3: manzana, pera, uva -> row with Stock=7, Precio=1.85).
# This shifts the "uva" row up to become row 3, matching the target layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Delete()
